$wb = $excel.ActiveWorkbook

$oldVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"
$newVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"

$aboutWs = $wb.Worksheets.Item("About")
$aboutWs.Range("A2").Value = "Version: $newVersion"
$aboutWs.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Vorgashorskaya Coal Mine, Russia, M0867, version '$newVersion'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

$dataWs = $wb.Worksheets.Item("Boundaries and methane sources")
for ($r = 2; $r -le 12; $r++) {
    $dataWs.Cells.Item($r, 19).Value = $newVersion
}
